# Update the "Authors" column (E) for rows 2-10 of Sheet1.
# The source data was re-ingested (elastic-search JSON export), which
# re-serialized these author lists with slightly different field padding;
# the cell text itself is refreshed to the newly-generated representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Fei%Zhou%NULL%0,                           Ting%Yu%NULL%0,                           Ronghui%Du%NULL%0,                           Guohui%Fan%NULL%0,                           Ying%Liu%NULL%0,                           Zhibo%Liu%NULL%0,                           Jie%Xiang%NULL%0,                           Yeming%Wang%NULL%0,                           Bin%Song%NULL%0,                           Xiaoying%Gu%NULL%0,                           Lulu%Guan%NULL%0,                           Yuan%Wei%NULL%0,                           Hui%Li%NULL%0,                           Xudong%Wu%NULL%0,                           Jiuyang%Xu%NULL%0,                           Shengjin%Tu%NULL%0,                           Yi%Zhang%NULL%0,                           Hua%Chen%NULL%0,                           Bin%Cao%NULL%0]'
$ws.Range("E3").Value = '[Litao%Zhang%NULL%1,                           Xinsheng%Yan%NULL%1,                           Qingkun%Fan%NULL%1,                           Haiyan%Liu%NULL%1,                           Xintian%Liu%NULL%1,                           Zejin%Liu%NULL%1,                           Zhenlu%Zhang%NULL%1]'
$ws.Range("E4").Value = '[Rong-Hui%Du%NULL%0,                           Li-Rong%Liang%NULL%0,                           Cheng-Qing%Yang%NULL%0,                           Wen%Wang%NULL%0,                           Tan-Ze%Cao%NULL%0,                           Ming%Li%NULL%0,                           Guang-Yun%Guo%NULL%0,                           Juan%Du%NULL%0,                           Chun-Lan%Zheng%NULL%0,                           Qi%Zhu%NULL%0,                           Ming%Hu%NULL%0,                           Xu-Yan%Li%NULL%0,                           Peng%Peng%NULL%0,                           Huan-Zhong%Shi%NULL%0]'
$ws.Range("E5").Value = '[Jianlei%Cao%NULL%0,                           Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%1,                           Wenlin%Cheng%NULL%1,                           Lei%Yu%NULL%1,                           Ya-Kun%Liu%NULL%1,                           Xiaoyong%Hu%NULL%1,                           Qiang%Liu%NULL%2]'
$ws.Range("E6").Value = '[Andrea%Giacomelli%andrea.giacomelli@unimi.it%0,                           Anna Lisa%Ridolfo%NULL%0,                           Laura%Milazzo%NULL%0,                           Letizia%Oreni%NULL%0,                           Dario%Bernacchia%NULL%0,                           Matteo%Siano%NULL%0,                           Cecilia%Bonazzetti%NULL%0,                           Alice%Covizzi%NULL%0,                           Marco%Schiuma%NULL%0,                           Matteo%Passerini%NULL%0,                           Marco%Piscaglia%NULL%0,                           Massimo%Coen%NULL%0,                           Guido%Gubertini%NULL%0,                           Giuliano%Rizzardini%NULL%0,                           Chiara%Cogliati%NULL%0,                           Anna Maria%Brambilla%NULL%0,                           Riccardo%Colombo%NULL%0,                           Antonio%Castelli%NULL%0,                           Roberto%Rech%NULL%0,                           Agostino%Riva%NULL%0,                           Alessandro%Torre%NULL%0,                           Luca%Meroni%NULL%0,                           Stefano%Rusconi%NULL%0,                           Spinello%Antinori%NULL%0,                           Massimo%Galli%NULL%0]'
$ws.Range("E7").Value = '[ I.%Paranjpe%null%1,                            A.% Russak%null%1,                            J. K.% De Freitas%null%1,                            A.% Lala%null%1,                            R.% Miotto%null%1,                            A.% Vaid%null%1,                            K. W.% Johnson%null%1,                            M.% Danieletto%null%1,                            E.% Golden%null%1,                            D.% Meyer%null%1,                            M.% Singh%null%1,                            S.% Somani%null%1,                            S.% Manna%null%1,                            U.% Nangia%null%1,                            A.% Kapoor%null%1,                            R.% O''Hagan%null%1,                            P. F.% O''Reilly%null%1,                            L. M.% Huckins%null%1,                            P.% Glowe%null%1,                            A.% Kia%null%1,                            P.% Timsina%null%1,                            R. M.% Freeman%null%1,                            M. A.% Levin%null%1,                            J.% Jhang%null%1,                            A.% Firpo%null%1,                            P.% Kovatch%null%1,                            J.% Finkelstein%null%1,                            J. A.% Aberg%null%1,                            E.% Bagiella%null%1,                            C. R.% Horowitz%null%1,                            B.% Murphy%null%1,                            Z. A.% Fayad%null%1,                            J.% Narula%null%1,                            E. J.% Nestler%null%1,                            V.% Fuster%null%1,                            C.% Cordon-Cardo%null%1,                            D. S.% Charney%null%1,                            D. L.% Reich%null%1,                            A. C.% Just%null%1,                            E. P.% Bottinger%null%1,                            A. W.% Charney%null%1,                            B. S.% Glicksberg%null%1,                            G.% Nadkarni%null%1,                           couldn''t find firstname% Mount Sina %null%1,         I.%Paranjpe%null%1,         A.% Russak%null%1,         J. K.% De Freitas%null%1,         A.% Lala%null%1,         R.% Miotto%null%1,         A.% Vaid%null%1,         K. W.% Johnson%null%1,         M.% Danieletto%null%1,         E.% Golden%null%1,         D.% Meyer%null%1,         M.% Singh%null%1,         S.% Somani%null%1,         S.% Manna%null%1,         U.% Nangia%null%1,         A.% Kapoor%null%1,         R.% O''Hagan%null%1,         P. F.% O''Reilly%null%1,         L. M.% Huckins%null%1,         P.% Glowe%null%1,         A.% Kia%null%1,         P.% Timsina%null%1,         R. M.% Freeman%null%1,         M. A.% Levin%null%1,         J.% Jhang%null%1,         A.% Firpo%null%1,         P.% Kovatch%null%0,         J.% Finkelstein%null%1,         J. A.% Aberg%null%1,         E.% Bagiella%null%1,         C. R.% Horowitz%null%1,         B.% Murphy%null%1,         Z. A.% Fayad%null%1,         J.% Narula%null%1,         E. J.% Nestler%null%1,         V.% Fuster%null%1,         C.% Cordon-Cardo%null%1,         D. S.% Charney%null%1,         D. L.% Reich%null%1,         A. C.% Just%null%1,         E. P.% Bottinger%null%1,         A. W.% Charney%null%1,         B. S.% Glicksberg%null%1,         G.% Nadkarni%null%1,        couldn''t find firstname% Mount Sina %null%1]'
$ws.Range("E8").Value = '[Tao%Chen%NULL%0,                           Di%Wu%NULL%5,                           Huilong%Chen%NULL%5,                           Weiming%Yan%NULL%5,                           Danlei%Yang%NULL%5,                           Guang%Chen%NULL%5,                           Ke%Ma%NULL%6,                           Dong%Xu%NULL%7,                           Haijing%Yu%NULL%5,                           Hongwu%Wang%NULL%5,                           Tao%Wang%NULL%0,                           Wei%Guo%NULL%5,                           Jia%Chen%NULL%5,                           Chen%Ding%NULL%5,                           Xiaoping%Zhang%NULL%5,                           Jiaquan%Huang%NULL%5,                           Meifang%Han%NULL%5,                           Shusheng%Li%NULL%5,                           Xiaoping%Luo%NULL%5,                           Jianping%Zhao%NULL%0,                           Qin%Ning%NULL%5]'
$ws.Range("E9").Value = '[ X.%Luo%null%1,                            H.% Xia%null%1,                            W.% Yang%null%1,                            B.% Wang%null%0,                            T.% Guo%null%1,                            J.% Xiong%null%1,                            Z.% Jiang%null%1,                            Y.% Liu%null%0,                            X.% Yan%null%1,                            W.% Zhou%null%1,                            L.% Ye%null%1,                            B. % Zhang%null%1,         X.%Luo%null%1,         H.% Xia%null%1,         W.% Yang%null%1,         B.% Wang%null%0,         T.% Guo%null%1,         J.% Xiong%null%1,         Z.% Jiang%null%1,         Y.% Liu%null%0,         X.% Yan%null%1,         W.% Zhou%null%0,         L.% Ye%null%1,         B. % Zhang%null%1]'
$ws.Range("E10").Value = '[Qingchun%Yao%xref no email%1,       Peng%Wang%xref no email%1,       Xingguang%Wang%xref no email%1,       Guoqiang%Qie%xref no email%1,       Mei%Meng%xref no email%1,       Xiwen%Tong%xref no email%1,       Xue%Bai%xref no email%1,       Min%Ding%xref no email%1,       Weiming%Liu%xref no email%1,       Keke%Liu%xref no email%1,       Yufeng%Chu%xref no email%1]'
